$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 235, shifting existing row 235 and below down by one.
$ws.Rows(235).Insert()

# Populate the newly inserted row 235 with the new data.
$ws.Cells.Item(235, 1).Value = 8
$ws.Cells.Item(235, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(235, 3).Value = "Coquimbo"
$ws.Cells.Item(235, 4).Value = 45229
$ws.Cells.Item(235, 5).Value = 4
$ws.Cells.Item(235, 6).Value = 100112001
$ws.Cells.Item(235, 7).Value = "Berenjena"
$ws.Cells.Item(235, 8).Value = "Sin especificar"
$ws.Cells.Item(235, 9).Value = "Primera"
$ws.Cells.Item(235, 10).Value = 360
$ws.Cells.Item(235, 11).Value = 11000
$ws.Cells.Item(235, 12).Value = 12000
$ws.Cells.Item(235, 13).Value = 11500
$ws.Cells.Item(235, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(235, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(235, 16).Value = 230
$ws.Cells.Item(235, 17).Value = 50
$ws.Cells.Item(235, 18).Value = "Hortaliza"
